# Append new data rows (A: 204..215, B: residual values) to Sheet1,
# extending the dataset from row 205 down to row 217.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A cells all share the same style (bold, bordered, centered).
# Copy the formatting from the last existing data cell (A205) onto the
# newly appended A cells so they match the rest of the column.
$styleSource = $ws.Range("A205")

$newRows = @(
    @{ Row = 206; A = 204; B = [double]"1.850371707708594E-17" },
    @{ Row = 207; A = 205; B = [double]"4.541821464375641E-17" },
    @{ Row = 208; A = 206; B = [double]"-6.938893903907228E-17" },
    @{ Row = 209; A = 207; B = [double]"6.167905692361981E-17" },
    @{ Row = 210; A = 208; B = [double]"0" },
    @{ Row = 211; A = 209; B = [double]"5.947623346206196E-17" },
    @{ Row = 212; A = 210; B = [double]"-4.625929269271486E-17" },
    @{ Row = 213; A = 211; B = [double]"-2.775557561562891E-17" },
    @{ Row = 214; A = 212; B = [double]"0" },
    @{ Row = 215; A = 213; B = [double]"0" },
    @{ Row = 216; A = 214; B = [double]"0" },
    @{ Row = 217; A = 215; B = [double]"0" }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $item.A
    $styleSource.Copy()
    $cellA.PasteSpecial(-4122)

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $item.B
}

$excel.CutCopyMode = $false
